$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    12354.73485800254,
    12333.06911207721,
    12326.04852480412,
    12302.48437123715,
    12276.4746913756,
    12256.6918218564,
    12242.03122744982,
    12227.57285863363,
    12221.50418426526,
    12212.26009819476,
    12196.41421434514,
    12181.11034186475
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
